$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: Postavljanje Expo projekta... -> Inicijalizacija I konfiguracija projekta ---
$ws.Range("B3").Value = "Inicijalizacija I konfiguracija projekta"
$ws.Range("C3").Value = "2h"
$ws.Range("D3").Value = "1h"
$ws.Range("E3").Value = "Leo Petrović"
$ws.Range("F3").Value = "Inicijaliziranje Expo projekta, stvaranje novog Supabase projekta, i spajanje na Supabase."

# --- Row 4: Implementacija ruta i osnovnog UI-a -> Implementacija navigacije i osnovnog UI-a ---
$ws.Range("B4").Value = "Implementacija navigacije i osnovnog UI-a"
$ws.Range("C4").Value = "4h"
$ws.Range("D4").Value = "5h"
$ws.Range("E4").Value = "Leo Petrović"
$ws.Range("F4").Value = "Konfigurirati rute za prijavu i registraciju, implementirati neke osnovne komponente, teme, boje, itd."

# --- Row 5: Implementacija prijave korisnika (name unchanged, note updated to Supabase) ---
$ws.Range("B5").Value = "Implementacija prijave korisnika"
$ws.Range("C5").Value = "3h"
$ws.Range("E5").Value = "Dragan Arapović"
$ws.Range("F5").Value = "Omogućiti prijavu korisnika pomoću Supabase-a."

# --- Row 6: Implementacija registracije korisnika ---
$ws.Range("B6").Value = "Implementacija registracije korisnika"
$ws.Range("C6").Value = "4h"
$ws.Range("E6").Value = "Mate Marić"
$ws.Range("F6").Value = "Dodati funkcionalnost registracije pomoću Supabase-a."

# --- Row 7: Postavljanje route guardova ---
$ws.Range("B7").Value = "Postavljanje route guardova"
$ws.Range("C7").Value = "3h"
$ws.Range("E7").Value = "Mate Marić"
$ws.Range("F7").Value = "Dodati zaštitu ruta za prijavljene korisnike koristeći Supabase realtime auth."

# --- Row 8: Testiranje funkcionalnosti prijave i registracije. ---
$ws.Range("B8").Value = "Testiranje funkcionalnosti prijave i registracije."
$ws.Range("C8").Value = "2h"
$ws.Range("E8").Value = "Dragan Arapović"
$ws.Range("F8").Value = "Provjeriti ispravnost ključnih funkcionalnosti."

# Rows 3, 4 and 7 now wrap onto a third line (longer "Napomena" text), so their
# row height grows from 25.5 (2 lines) to 38.25 (3 lines) to match the re-wrapped content.
$ws.Rows.Item(3).RowHeight = 38.25
$ws.Rows.Item(4).RowHeight = 38.25
$ws.Rows.Item(7).RowHeight = 38.25

# Selection cursor moved to H3 on the first sheet.
$ws.Range("H3").Select()
